$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (column A values are Excel serial date numbers,
# same representation/format as the existing rows above them)
$data = @(
    @(367, 44441, 1, 21, 116.809433752364),
    @(368, 44442, 5, 26, 144.621203693403),
    @(369, 44443, 2, 24, 133.4964957169874),
    @(370, 44444, 3, 19, 105.6847257759484),
    @(371, 44445, 5, 19, 105.6847257759484),
    @(372, 44446, 1, 17, 94.56001779953276),
    @(373, 44447, 0, 17, 94.56001779953276),
    @(374, 44448, 15, 31, 172.4329736344421)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy the formatting of column A from the row above (date style, font, border, alignment)
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = 0
